# Apply updated cryptos list values (Price + Volume(1h)) to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.078.34"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.396.06"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.34%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.63"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.31"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.06%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.84%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.395.25"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.413"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.988.82"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.51%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.01"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.148.81"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.407.49"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.87"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "366.34"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.53"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.86"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.69%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.74"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.179"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.98"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.16"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.14%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.61"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.81%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.13"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -8.26%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.683.93"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.34"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.26"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0676"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "335.33"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +9.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.69"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.48"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.40"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.96%  "
